$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Article sur l économie"
$ws.Range("C2").Value = "L économie mondiale montre des signes de reprise..."

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-01-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Actualités politiques"
$ws.Range("C3").Value = "Les dernières décisions politiques impactent..."

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2024-02-01"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "Innovations technologiques"
$ws.Range("C4").Value = "Les nouvelles technologies révolutionnent..."

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-02-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "Développement durable"
$ws.Range("C5").Value = "Les initiatives environnementales progressent..."
